$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for "Choclo" (Choclero, Primera)
# at Terminal Hortofrutícola Agro Chillán. It belongs right after the
# existing row 277 (the "Dulce o Americano" entry), so insert a new row at
# 278 - this pushes the old rows 278-303 down to 279-304 (matches the diff).
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new record. The columns that
# are constant for every "Choclo" row at this market (A, B, C, E, F, G, H,
# N, Q, R) are carried over unchanged; only the date, quality, prices and
# origin (D, I, J, K, L, M, O, P) are new values for this week.
$ws.Cells.Item(278, 1).Value = 7
$ws.Cells.Item(278, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(278, 3).Value = "Ñuble"
$ws.Cells.Item(278, 4).Value = 45013
$ws.Cells.Item(278, 5).Value = 16
$ws.Cells.Item(278, 6).Value = 100112024
$ws.Cells.Item(278, 7).Value = "Choclo"
$ws.Cells.Item(278, 8).Value = "Choclero"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 10000
$ws.Cells.Item(278, 11).Value = 400
$ws.Cells.Item(278, 12).Value = 450
$ws.Cells.Item(278, 13).Value = 430
$ws.Cells.Item(278, 14).Value = "`$/unidad"
$ws.Cells.Item(278, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(278, 16).Value = 430
$ws.Cells.Item(278, 17).Value = 1
$ws.Cells.Item(278, 18).Value = "Hortaliza"
